# Arrange files in order and commenting code
#
# The source data list (rows 2-28 of Sheet1) is re-sorted/grouped so that all
# dosage-strength variants of the same brand sit together, the now-unused
# "Levomax" brand/row is removed, and the final blank/duplicate row (row 29)
# is dropped. Shared strings that become orphaned (e.g. "Levomax",
# "Levomax 750mg Tablet - 10's", "10 's") are automatically compacted out of
# the workbook once nothing references them anymore.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Re-order the "Dinafex" tablet rows (swap 180mg/120mg) ---
$ws.Range("D3").Value = "Dinafex 120mg Tablet"
$ws.Range("D4").Value = "Dinafex 180mg Tablet"

# --- Re-order the "Etorix" tablet rows (90mg now before 120mg) ---
$ws.Range("D7").Value = "Etorix 90mg Tablet"
$ws.Range("E7").Value = "30's"
$ws.Range("D9").Value = "Etorix 120mg Tablet"
$ws.Range("E9").Value = "20's"

# --- Re-order the "Flucloxin" capsule rows (36's pack now first) ---
$ws.Range("D11").Value = "Flucloxin 500mg Capsule - 36's"
$ws.Range("E11").Value = "36 's"
$ws.Range("D12").Value = "Flucloxin 500mg Capsule"
$ws.Range("E12").Value = "30 's"

# --- Re-order the "Ketonic" rows (Injection, IM/IV Injection, Tablet) ---
$ws.Range("D14").Value = "Ketonic 30mg Injection"
$ws.Range("E14").Value = "5 's"
$ws.Range("D15").Value = "Ketonic 30mg IM/IV Injection - 4's"
$ws.Range("E15").Value = "4's"
$ws.Range("D16").Value = "Ketonic 10mg Tablet"
$ws.Range("E16").Value = "20's"

# --- Drop the "Levomax" brand row entirely; shift the remaining brands
#     (Naprox, Oradin, Osticare, Rupaday, Sk-Mox, Zithrox) up by one row ---
$ws.Range("A20").Value = 17
$ws.Range("B20").Value = "Naprox"
$ws.Range("D20").Value = "Naprox Plus 500mg Tablet - 30's"
$ws.Range("E20").Value = "30 's"

$ws.Range("A21").Value = 19
$ws.Range("B21").Value = "Oradin"
$ws.Range("D21").Value = "Oradin Plus Tablet - 40's"
$ws.Range("E21").Value = "40 's"

$ws.Range("A22").Value = 20
$ws.Range("B22").Value = "Osticare"
$ws.Range("D22").Value = "Osticare Tablet 24's"
$ws.Range("E22").Value = "24's"

$ws.Range("A23").Value = 23
$ws.Range("B23").Value = "Rupaday"
$ws.Range("D23").Value = "Rupaday Oral Solution 60ml"
$ws.Range("E23").Value = "1's"

$ws.Range("A24").Value = 24
$ws.Range("B24").Value = "Sk-Mox"
$ws.Range("D24").Value = "Sk-Mox 500mg Capsule"
$ws.Range("E24").Value = "48 's"

$ws.Range("A25").Value = 35
$ws.Range("B25").Value = "Zithrox"
$ws.Range("D25").Value = "Zithrox 250mg Tablet - 6's"
$ws.Range("E25").Value = "6's"

# --- Re-order the remaining "Zithrox" rows (500mg tablet, 15ml, 30ml) ---
$ws.Range("D26").Value = "Zithrox 500mg Tablet"
$ws.Range("E26").Value = "6 's"
$ws.Range("D27").Value = "Zithrox 15ml Suspension"
$ws.Range("E27").Value = "15 ml"
$ws.Range("D28").Value = "Zithrox 30ml Dry Suspension"
$ws.Range("E28").Value = "30ml"

# --- Remove the now-superfluous last row (old row 29) ---
$ws.Rows(29).Delete() | Out-Null
